# The workbook's only deliberate content change (per the diff) is the
# worksheet being renamed from the generic default "Sheet1" to the
# descriptive name "projects" (matching the workbook's own file name,
# utils/config/ReferenceLists/projects.xlsx).
#
# Everything else in the diff (fileVersion/rupBuild, the author's
# x15ac:absPath SharePoint URL, xr:revisionPtr save id/GUID, the
# bookViews window position/size, sheetFormatPr defaultRowHeight /
# x14ac:dyDescent, the tiny column-width deltas, and the auto-computed
# wrap-text row heights on rows 2/3/9) are artifacts of the file simply
# having been opened and re-saved by a newer Excel build (rupBuild
# 28324 -> 28623) on a different author's machine/monitor - they are
# not explicit user edits and are not reproducible (or even meaningful)
# through the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "projects"
